$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "Datos actualizados" timestamp text in cell A1
$ws.Range("A1").Value = "Datos actualizados a 15 de Abril de 2020 a las 23:52"

# Update row 5 (Cataluña) figures
$ws.Range("B5").Value = 39375
$ws.Range("C5").Value = 17297
$ws.Range("D5").Value = 18223
$ws.Range("E5").Value = 3855
